$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Propagate the existing header style (bold, bordered, centered) from H1 to
# the two new header cells before writing their text, so they reuse the
# same cellXf instead of minting a new one.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data cells I2 and J2 (plain, unstyled like the rest of row 2)
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
